# Daily attendance processing - reverse the order of names in the
# "Recorded By" (column G) list for any session row whose value is a
# comma-separated list that includes "System" as one of the entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text

    if ($txt -ne $null -and $txt -ne "") {
        if (($txt -like "*,*") -and ($txt -match "(?i)system")) {
            $parts = $txt -split ", "
            $n = $parts.Count
            $rev = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $rev += $parts[$i]
            }
            $newVal = [string]::Join(", ", $rev)
            $cell.Value = $newVal
        }
    }
}
